# Update "想去人数" (interested-count) figures in the "展览" and "全部类型" sheets.
# Mapping of F-column (row -> new value) per sheet, derived from the source diff.

$wb = $excel.ActiveWorkbook

$sheet1Updates = @{
    2  = 4963
    3  = 105
    4  = 345
    5  = 45
    6  = 18
    7  = 50
    8  = 129
    10 = 320
    11 = 260
    12 = 2979
    13 = 156
    14 = 1585
    15 = 12
}

$sheet4Updates = @{
    2  = 4963
    3  = 105
    4  = 345
    5  = 45
    7  = 18
    8  = 50
    9  = 129
    11 = 320
    12 = 260
    13 = 2979
    14 = 156
    15 = 1585
    16 = 12
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
